# Replace the Michaelis-Menten model's computed (formula-driven) rows with a
# fixed table of experimentally measured concentration/velocity pairs, and
# drop the now-unused extra rows (15-17) so the sheet only spans A1:F14.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New concentration (A) / measured velocity (B) values for rows 3-14.
# These become plain literal values (the old shared formula in B3:B17 is
# overwritten / removed).
$data = @(
    @{ Row = 3;  A = 0.05; B = 1.446 }
    @{ Row = 4;  A = 0.1;  B = 2.6469999999999998 }
    @{ Row = 5;  A = 0.2;  B = 4.6719999999999997 }
    @{ Row = 6;  A = 0.3;  B = 7.7729999999999997 }
    @{ Row = 7;  A = 0.5;  B = 11.432 }
    @{ Row = 8;  A = 1;    B = 18.792000000000002 }
    @{ Row = 9;  A = 1.5;  B = 23.913 }
    @{ Row = 10; A = 2;    B = 26.5505 }
    @{ Row = 11; A = 3;    B = 29.311499999999999 }
    @{ Row = 12; A = 4;    B = 32.101500000000001 }
    @{ Row = 13; A = 5;    B = 31.4345 }
    @{ Row = 14; A = 7.5;  B = 37.161999999999999 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.A
    $ws.Cells.Item($r, 2).Value = $entry.B
}

# Rows 15-17 no longer exist in the new layout - remove them entirely.
$ws.Range("A15:F17").EntireRow.Delete()

# Update the selection to match the saved state (B14).
$ws.Range("B14").Select()
